$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - keep existing style (s="1"), just set values
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "LinkedIn_Poster"
$ws.Range("G1").Value = "LinkedIn_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

# Data row (row 2)
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "QA Analyst"
$ws.Range("C2").Value = "A QA Analyst job description involves creating, executing, and maintaining test plans and cases to ensure software quality, collaborating with cross-functional teams to resolve issues, mentoring junior analysts, and improving overall quality processes. Key responsibilities include designing test strategies, automating tests, reporting bugs, performing root cause analysis, and ensuring quality throughout the development lifecycle. "
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 3
